# Update the "想去人数" (F column) counts that changed between the two
# data-refresh snapshots of the generated gh-pages output.
#
# The workbook has two sheets that carry the same underlying rows:
#   - "展览"    (sheet1 / rId1)
#   - "全部类型" (sheet4 / rId4)
# Both need the same F-column value bumps, just on different row numbers.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> list of (cell address, new value)
$updates = @{
    "展览" = @(
        @{ Cell = "F4";  Value = 61 }
        @{ Cell = "F5";  Value = 1639 }
        @{ Cell = "F7";  Value = 793 }
        @{ Cell = "F8";  Value = 2014 }
        @{ Cell = "F9";  Value = 1936 }
        @{ Cell = "F10"; Value = 989 }
        @{ Cell = "F11"; Value = 346 }
        @{ Cell = "F13"; Value = 1596 }
        @{ Cell = "F14"; Value = 342 }
        @{ Cell = "F18"; Value = 1412 }
        @{ Cell = "F19"; Value = 506 }
        @{ Cell = "F21"; Value = 310 }
        @{ Cell = "F22"; Value = 10578 }
        @{ Cell = "F23"; Value = 9746 }
        @{ Cell = "F26"; Value = 1819 }
    )
    "全部类型" = @(
        @{ Cell = "F6";  Value = 61 }
        @{ Cell = "F7";  Value = 1639 }
        @{ Cell = "F9";  Value = 793 }
        @{ Cell = "F10"; Value = 2014 }
        @{ Cell = "F11"; Value = 1936 }
        @{ Cell = "F12"; Value = 989 }
        @{ Cell = "F13"; Value = 346 }
        @{ Cell = "F15"; Value = 1596 }
        @{ Cell = "F16"; Value = 342 }
        @{ Cell = "F22"; Value = 1412 }
        @{ Cell = "F23"; Value = 506 }
        @{ Cell = "F25"; Value = 310 }
        @{ Cell = "F26"; Value = 10578 }
        @{ Cell = "F27"; Value = 9746 }
        @{ Cell = "F30"; Value = 1819 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Range($entry.Cell).Value = $entry.Value
    }
}
